# Update the LSTM baseline results row (row 2) with the re-run metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 31                     # Elementos
$ws.Range("G2").Value = "adam"                 # Optimizador (was the raw Keras repr)
$ws.Range("I2").Value = 32                     # BatchSize
$ws.Range("J2").Value = 47.56829888765251      # MAE [$COP/kWh]
$ws.Range("K2").Value = 3747.00177240947       # MSE [$COP/kWh]
$ws.Range("L2").Value = 61.21275824866471      # RMSE [$COP/kWh]
$ws.Range("M2").Value = 0.2873282370544794     # MAPE [%]
